$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep/become a text value (avoids Excel
    # auto-converting numeric-looking strings like "1.00" or "14.00"
    # into numbers), while preserving the cell's original style index.
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '46.896.89'
Set-TextValue $ws.Range('E2') '  +6.81%  '
Set-TextValue $ws.Range('D3') '2.309.25'
Set-TextValue $ws.Range('E3') '  +5.31%  '
Set-TextValue $ws.Range('E4') '  -0.46%  '
Set-TextValue $ws.Range('D5') '298.37'
Set-TextValue $ws.Range('E5') '  +1.61%  '
Set-TextValue $ws.Range('D6') '100.73'
Set-TextValue $ws.Range('E6') '  +15.02%  '
Set-TextValue $ws.Range('E7') '  +1.56%  '
Set-TextValue $ws.Range('D8') '1.00'
Set-TextValue $ws.Range('E8') '  -0.46%  '
Set-TextValue $ws.Range('E9') '  +11.04%  '
Set-TextValue $ws.Range('D10') '36.21'
Set-TextValue $ws.Range('E10') '  +13.74%  '
Set-TextValue $ws.Range('D11') '0.0800'
Set-TextValue $ws.Range('E11') '  +5.04%  '
Set-TextValue $ws.Range('D12') '7.34'
Set-TextValue $ws.Range('E12') '  +9.72%  '
Set-TextValue $ws.Range('E13') '  +1.72%  '
Set-TextValue $ws.Range('D14') '2.661.15'
Set-TextValue $ws.Range('E14') '  +5.31%  '
Set-TextValue $ws.Range('D15') '2.307.08'
Set-TextValue $ws.Range('E15') '  +2.43%  '
Set-TextValue $ws.Range('D16') '14.00'
Set-TextValue $ws.Range('E16') '  +8.91%  '
Set-TextValue $ws.Range('D17') '0.819'
Set-TextValue $ws.Range('E17') '  +7.39%  '
Set-TextValue $ws.Range('D18') '46.848.34'
Set-TextValue $ws.Range('E18') '  +7.79%  '
Set-TextValue $ws.Range('D19') '13.16'
Set-TextValue $ws.Range('E19') '  +24.66%  '
Set-TextValue $ws.Range('D20') '0.0₃0942'
Set-TextValue $ws.Range('E20') '  +7.43%  '
Set-TextValue $ws.Range('D21') '6.16'
Set-TextValue $ws.Range('E21') '  +6.58%  '
Set-TextValue $ws.Range('D22') '66.99'
Set-TextValue $ws.Range('E22') '  +7.11%  '
Set-TextValue $ws.Range('D23') '248.47'
Set-TextValue $ws.Range('E23') '  +8.48%  '
Set-TextValue $ws.Range('E24') '  +6.68%  '
Set-TextValue $ws.Range('D25') '1.98'
Set-TextValue $ws.Range('E25') '  +10.17%  '
Set-TextValue $ws.Range('E26') '  -0.36%  '
Set-TextValue $ws.Range('D27') '42.91'
Set-TextValue $ws.Range('E27') '  +21.96%  '
Set-TextValue $ws.Range('E28') '  +1.71%  '
Set-TextValue $ws.Range('D29') '9.92'
Set-TextValue $ws.Range('E29') '  +8.58%  '
Set-TextValue $ws.Range('D30') '20.26'
Set-TextValue $ws.Range('E30') '  +6.50%  '
Set-TextValue $ws.Range('D31') '5.77'
Set-TextValue $ws.Range('E31') '  +9.99%  '
Set-TextValue $ws.Range('D32') '147.05'
Set-TextValue $ws.Range('E32') '  +0.71%  '
Set-TextValue $ws.Range('D33') '0.0800'
Set-TextValue $ws.Range('E33') '  +10.47%  '
Set-TextValue $ws.Range('E34') '  +4.52%  '
Set-TextValue $ws.Range('E35') '  +11.54%  '
Set-TextValue $ws.Range('D36') '3.11'
Set-TextValue $ws.Range('E36') '  +9.32%  '
Set-TextValue $ws.Range('E37') '  +3.52%  '
Set-TextValue $ws.Range('D38') '1.79'
Set-TextValue $ws.Range('E38') '  +9.99%  '
Set-TextValue $ws.Range('D39') '15.74'
Set-TextValue $ws.Range('E39') '  +21.73%  '
Set-TextValue $ws.Range('E40') '  +16.36%  '
Set-TextValue $ws.Range('D41') '3.40'
Set-TextValue $ws.Range('E41') '  +12.47%  '
Set-TextValue $ws.Range('E42') '  +10.84%  '
Set-TextValue $ws.Range('E44') '  +22.64%  '
Set-TextValue $ws.Range('D45') '1.839.53'
Set-TextValue $ws.Range('E45') '  +4.87%  '
Set-TextValue $ws.Range('D46') '90.08'
Set-TextValue $ws.Range('E46') '  +24.60%  '
Set-TextValue $ws.Range('D47') '0.199'
Set-TextValue $ws.Range('E47') '  +16.82%  '
Set-TextValue $ws.Range('D48') '75.48'
Set-TextValue $ws.Range('E48') '  +15.78%  '
Set-TextValue $ws.Range('E49') '  +10.97%  '
Set-TextValue $ws.Range('D50') '97.32'
Set-TextValue $ws.Range('E50') '  +6.93%  '
Set-TextValue $ws.Range('D51') '54.25'
Set-TextValue $ws.Range('E51') '  +12.61%  '
